$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, copying the style from the existing header (H1)
# so that the new header cells match the existing bold/bordered/centered style.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I2:J67 data block with the new I0/IF values for each row.
$arr = New-Object 'object[,]' 66,2
$arr[0,0] = 8; $arr[0,1] = 8
$arr[1,0] = 5; $arr[1,1] = 5
$arr[2,0] = 9; $arr[2,1] = 9
$arr[3,0] = 7; $arr[3,1] = 7
$arr[4,0] = 8; $arr[4,1] = 8
$arr[5,0] = 5; $arr[5,1] = 5
$arr[6,0] = 8; $arr[6,1] = 8
$arr[7,0] = 5; $arr[7,1] = 5
$arr[8,0] = 6; $arr[8,1] = 6
$arr[9,0] = 8; $arr[9,1] = 8
$arr[10,0] = 7; $arr[10,1] = 7
$arr[11,0] = 2; $arr[11,1] = 3
$arr[12,0] = 8; $arr[12,1] = 8
$arr[13,0] = 3; $arr[13,1] = 3
$arr[14,0] = 9; $arr[14,1] = 9
$arr[15,0] = 8; $arr[15,1] = 8
$arr[16,0] = 8; $arr[16,1] = 8
$arr[17,0] = 6; $arr[17,1] = 7
$arr[18,0] = 8; $arr[18,1] = 9
$arr[19,0] = 6; $arr[19,1] = 6
$arr[20,0] = 7; $arr[20,1] = 8
$arr[21,0] = 8; $arr[21,1] = 8
$arr[22,0] = 8; $arr[22,1] = 9
$arr[23,0] = 6; $arr[23,1] = 6
$arr[24,0] = 2; $arr[24,1] = 3
$arr[25,0] = 6; $arr[25,1] = 6
$arr[26,0] = 3; $arr[26,1] = 4
$arr[27,0] = 7; $arr[27,1] = 8
$arr[28,0] = 8; $arr[28,1] = 8
$arr[29,0] = 9; $arr[29,1] = 9
$arr[30,0] = 7; $arr[30,1] = 7
$arr[31,0] = 5; $arr[31,1] = 5
$arr[32,0] = 6; $arr[32,1] = 7
$arr[33,0] = 6; $arr[33,1] = 6
$arr[34,0] = 4; $arr[34,1] = 4
$arr[35,0] = 5; $arr[35,1] = 6
$arr[36,0] = 6; $arr[36,1] = 6
$arr[37,0] = 6; $arr[37,1] = 6
$arr[38,0] = 6; $arr[38,1] = 6
$arr[39,0] = 7; $arr[39,1] = 7
$arr[40,0] = 6; $arr[40,1] = 6
$arr[41,0] = 6; $arr[41,1] = 7
$arr[42,0] = 8; $arr[42,1] = 8
$arr[43,0] = 9; $arr[43,1] = 9
$arr[44,0] = 8; $arr[44,1] = 8
$arr[45,0] = 4; $arr[45,1] = 4
$arr[46,0] = 5; $arr[46,1] = 5
$arr[47,0] = 7; $arr[47,1] = 7
$arr[48,0] = 5; $arr[48,1] = 5
$arr[49,0] = 5; $arr[49,1] = 5
$arr[50,0] = 6; $arr[50,1] = 6
$arr[51,0] = 8; $arr[51,1] = 8
$arr[52,0] = 6; $arr[52,1] = 6
$arr[53,0] = 4; $arr[53,1] = 5
$arr[54,0] = 7; $arr[54,1] = 8
$arr[55,0] = 8; $arr[55,1] = 8
$arr[56,0] = 4; $arr[56,1] = 5
$arr[57,0] = 8; $arr[57,1] = 8
$arr[58,0] = 6; $arr[58,1] = 6
$arr[59,0] = 4; $arr[59,1] = 4
$arr[60,0] = 5; $arr[60,1] = 5
$arr[61,0] = 6; $arr[61,1] = 6
$arr[62,0] = 7; $arr[62,1] = 7
$arr[63,0] = 5; $arr[63,1] = 5
$arr[64,0] = 3; $arr[64,1] = 3
$arr[65,0] = 4; $arr[65,1] = 4
$ws.Range("I2:J67").Value = $arr
